# Applies the "Added more accurate values" commit:
#  - Param!G5 changes from 10 to 1
#  - Each result sheet ("1","2","3","4") gains a new column C "Li-Wu" holding
#    the previously-displayed (rounded) price, while column B is updated to
#    hold a more precise value. Column widths / number formats follow.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Param sheet: G5 10 -> 1
# ---------------------------------------------------------------------
$wsParam = $wb.Worksheets.Item("Param")
$wsParam.Activate()
$wsParam.Range("G5").Value = 1

# ---------------------------------------------------------------------
# Helper data for the four result sheets: sheet name -> rows of
# (old displayed value, new more-precise value)
# ---------------------------------------------------------------------

# Sheet "1"
$ws1 = $wb.Worksheets.Item("1")
$ws1.Range("C1").Value = "Li-Wu"
$ws1.Range("C2").Value = 13.214930000000001
$ws1.Range("C2").NumberFormat = "0.00000"
$ws1.Range("B2").Value = 13.214919999999999
$ws1.Range("B2").NumberFormat = "0.000000"
$ws1.Columns.Item(2).ColumnWidth = 10.90625
$ws1.Columns.Item(3).ColumnWidth = 9.6328125

# Sheet "2"
$ws2 = $wb.Worksheets.Item("2")
$ws2.Range("C1").Value = "Li-Wu"
$ws2.Range("C2").Value = 40.797730000000001
$ws2.Range("C2").NumberFormat = "0.00000"
$ws2.Range("B2").Value = 40.797688999999998
$ws2.Range("B2").NumberFormat = "0.000000"
$ws2.Columns.Item(2).ColumnWidth = 10.90625
$ws2.Columns.Item(3).ColumnWidth = 9.6328125

# Sheet "3"
$ws3 = $wb.Worksheets.Item("3")
$ws3.Range("C1").Value = "Li-Wu"
$ws3.Range("C2").Value = 62.763120000000001
$ws3.Range("C2").NumberFormat = "0.00000"
$ws3.Range("B2").Value = 62.763120000000001
$ws3.Range("B2").NumberFormat = "0.000000"
$ws3.Columns.Item(2).ColumnWidth = 10.90625
$ws3.Columns.Item(3).ColumnWidth = 9.6328125
$ws3.Activate()
$ws3.Columns.Item(3).Select()

# Sheet "4"
$ws4 = $wb.Worksheets.Item("4")
$ws4.Range("C1").Value = "Li-Wu"

$ws4.Range("C2").Value = 21.41873
$ws4.Range("C2").NumberFormat = "0.00000"
$ws4.Range("B2").Value = 21.418717000000001
$ws4.Range("B2").NumberFormat = "0.000000"

$ws4.Range("C3").Value = 15.16798
$ws4.Range("B3").Value = 15.16797
$ws4.Range("B3").NumberFormat = "0.000000"

$ws4.Range("C4").Value = 10.174480000000001
$ws4.Range("B4").Value = 10.174469
$ws4.Range("B4").NumberFormat = "0.000000"

$ws4.Columns.Item(2).ColumnWidth = 10.90625
$ws4.Columns.Item(3).ColumnWidth = 9.6328125
$ws4.Activate()
$ws4.Range("B2").Select()

# ---------------------------------------------------------------------
# Leave "Param" as the active sheet, matching the original workbook state.
# ---------------------------------------------------------------------
$wsParam.Activate()
